# Updated CVDs for the month
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rosemont Illinois")

# Row 4 updates (columns O through W)
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# Row 7 update: O7 becomes blank
$ws.Range("O7").ClearContents()
